# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with latest feed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.585.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "'2.260.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'250.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'75.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.88%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.638"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("D10").Value = "'39.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "'0.0969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'7.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'2.598.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'14.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'0.861"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "'2.256.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'42.492.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "'72.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "'235.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'2.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'3.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.87%  "
$ws.Range("D26").Value = "'11.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "'167.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'20.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "'6.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").Value = "'0.0854"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.57%  "
$ws.Range("D33").Value = "'0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "'31.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'4.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("D39").Value = "'13.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.18%  "
$ws.Range("E40").Value = "  -3.42%  "
$ws.Range("D41").Value = "'5.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'0.206"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'61.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("D44").Value = "'107.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.49%  "
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").Value = "'4.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.22%  "
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "'1.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("E51").Value = "  -2.39%  "
